$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to remain plain text so values such as "1.00" or
# "604.50" are not silently normalised into numbers (losing trailing
# zeros / dropping the decimal grouping dots used by the source site).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.690.18"
$ws.Range("E2").Value = "  +2.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.204.74"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.50"
$ws.Range("E5").Value = "  +4.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.03"
$ws.Range("E6").Value = "  +4.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.559"
$ws.Range("E8").Value = "  +6.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.204.46"
$ws.Range("E9").Value = "  +1.52%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.88"
$ws.Range("E11").Value = "  -4.39%  "

$ws.Range("E12").Value = "  +3.89%  "

$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.37"
$ws.Range("E14").Value = "  +5.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.729.50"
$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.675.94"
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("E17").Value = "  +5.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.204.47"
$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "525.71"
$ws.Range("E19").Value = "  +3.77%  "

$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.57"
$ws.Range("E21").Value = "  +4.23%  "

$ws.Range("E22").Value = "  +3.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.25"
$ws.Range("E23").Value = "  +6.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.07"
$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.77"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  +2.28%  "

$ws.Range("E28").Value = "  +3.03%  "

$ws.Range("E29").Value = "  +9.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  +6.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").Value = "  +8.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.38"
$ws.Range("E32").Value = "  +2.66%  "

$ws.Range("E33").Value = "  +3.12%  "

$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "517.88"
$ws.Range("E36").Value = "  +7.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.92"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0909"
$ws.Range("E38").Value = "  +2.20%  "

$ws.Range("E39").Value = "  +2.70%  "

$ws.Range("E40").Value = "  +8.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.93"
$ws.Range("E41").Value = "  +2.19%  "

$ws.Range("E42").Value = "  -1.20%  "

$ws.Range("E43").Value = "  +14.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.302"
$ws.Range("E44").Value = "  +6.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.899.21"
$ws.Range("E46").Value = "  -3.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.71"
$ws.Range("E47").Value = "  +1.08%  "

$ws.Range("E48").Value = "  +11.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.119"
$ws.Range("E49").Value = "  +3.88%  "

$ws.Range("E51").Value = "  +4.59%  "
